# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (same header/format template as the
#    new quarter) into a new sheet positioned right before "总计",
#    rename it "2022-Q1", and overwrite its holdings with the 2022-Q1
#    figures.
# 2. Insert a new top row into "总计" for the "2022-Q1" summary line and
#    renumber the existing rows' running index.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) New "2022-Q1" worksheet
# ---------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($totalSheet)

$ws = $wb.Worksheets.Item("2021-Q4 (2)")
$ws.Name = "2022-Q1"

# Force text storage (keeps leading zeros / fixed decimal strings) for
# the fund-code column and the numeric-looking text columns, matching
# how the source data is published (everything except the rank is text).
$ws.Range("B2:B4").NumberFormat = "@"
$ws.Range("D2:G4").NumberFormat = "@"

$ws.Range("B2").Value = "002423"
$ws.Range("C2").Value = "华宝兴业标普美国消费(QDII-LOF)美元"
$ws.Range("D2").Value = "3.62"
$ws.Range("E2").Value = "94.37"
$ws.Range("F2").Value = "2.39"
$ws.Range("G2").Value = "0.0865"
$ws.Range("H2").Value = 7

$ws.Range("B3").Value = "162415"
$ws.Range("C3").Value = "华宝标普美国消费(QDII-LOF)人民币A"
$ws.Range("D3").Value = "3.62"
$ws.Range("E3").Value = "94.37"
$ws.Range("F3").Value = "2.39"
$ws.Range("G3").Value = "0.0865"
$ws.Range("H3").Value = 7

$ws.Range("B4").Value = "009975"
$ws.Range("C4").Value = "华宝标普美国消费(QDII-LOF)人民币C"
$ws.Range("D4").Value = "0.61"
$ws.Range("E4").Value = "94.37"
$ws.Range("F4").Value = "2.39"
$ws.Range("G4").Value = "0.0146"
$ws.Range("H4").Value = 7

# ---------------------------------------------------------------
# 2) Update the "总计" summary sheet with the new quarter's row
# ---------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()
$total.Range("A2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.19

# Re-apply the running-index column's style (copied from the row below,
# which already carries the correct formatting for column A).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

# Renumber the running index of the rows that got pushed down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# Restore the originally active sheet/tab.
$wb.Worksheets.Item("2020-Q4").Activate()
